$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its existing text formatting so numeric-looking
# price strings (e.g. "1.001", "0.02350") are not silently converted to
# actual numbers (which would normalize/round them).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.600.09'
$ws.Range("D3").Value = '1.829.80'
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '316.34'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = '0.5332'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.3988'
$ws.Range("E8").Value = '  +5.39%  '
$ws.Range("D9").Value = '0.07829'
$ws.Range("E9").Value = '  +4.60%  '
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").Value = '1.119'
$ws.Range("E11").Value = '  +1.97%  '
$ws.Range("D12").Value = '6.348'
$ws.Range("E12").Value = '  +2.02%  '
$ws.Range("D13").Value = '21.07'
$ws.Range("E13").Value = '  +2.48%  '
$ws.Range("D14").Value = '7.591'
$ws.Range("E14").Value = '  +2.83%  '
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("D16").Value = '1.828.70'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("D17").Value = '93.38'
$ws.Range("E17").Value = '  +3.97%  '
$ws.Range("D18").Value = '0.00001095'
$ws.Range("E18").Value = '  +2.97%  '
$ws.Range("D19").Value = '0.06564'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = '17.83'
$ws.Range("E20").Value = '  +2.20%  '
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").Value = '6.109'
$ws.Range("E22").Value = '  +3.06%  '
$ws.Range("D23").Value = '28.606.99'
$ws.Range("D24").Value = '11.24'
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("E25").Value = '  +7.22%  '
$ws.Range("D26").Value = '20.86'
$ws.Range("E26").Value = '  +1.67%  '
$ws.Range("D27").Value = '157.16'
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").Value = '2.040.52'
$ws.Range("E28").Value = '  +1.22%  '
$ws.Range("D29").Value = '2.414'
$ws.Range("E29").Value = '  +3.86%  '
$ws.Range("D30").Value = '125.55'
$ws.Range("E30").Value = '  +2.82%  '
$ws.Range("D31").Value = '1.147'
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("D32").Value = '0.1122'
$ws.Range("E32").Value = '  +2.34%  '
$ws.Range("D33").Value = '5.744'
$ws.Range("E33").Value = '  +2.73%  '
$ws.Range("D34").Value = '3.653'
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("D35").Value = '0.07305'
$ws.Range("E35").Value = '  +1.95%  '
$ws.Range("D36").Value = '0.2269'
$ws.Range("E36").Value = '  +1.80%  '
$ws.Range("D37").Value = '8.989'
$ws.Range("E37").Value = '  +6.39%  '
$ws.Range("D38").Value = '0.02350'
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("D39").Value = '5.228'
$ws.Range("E39").Value = '  +2.48%  '
$ws.Range("D40").Value = '11.41'
$ws.Range("E40").Value = '  +2.48%  '
$ws.Range("D41").Value = '0.6304'
$ws.Range("E41").Value = '  +2.03%  '
$ws.Range("D42").Value = '1.197'
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("E44").Value = '  -2.82%  '
$ws.Range("D45").Value = '13.57'
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("D46").Value = '0.5940'
$ws.Range("E46").Value = '  +2.88%  '
$ws.Range("D47").Value = '3.715'
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").Value = '125.57'
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("D49").Value = '1.999'
$ws.Range("E49").Value = '  +3.66%  '
$ws.Range("D50").Value = '1.193'
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("D51").Value = '0.06953'
$ws.Range("E51").Value = '  +1.96%  '
